$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 155.8700226666667
$ws.Range("H2").Value = 467.610068
$ws.Range("I2").Value = 0.4627663557222626
$ws.Range("J2").Value = 0.4864916976605717
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1613523333333333
$ws.Range("N2").Value = 0.484057
$ws.Range("O2").Value = 0.2926766298022186
$ws.Range("P2").Value = 0.3782017374917083
$ws.Range("Q2").Value = 25.14999185398623
$ws.Range("R2").Value = 226.349926685876
$ws.Range("S2").Value = 0.1354408973786465
$ws.Range("T2").Value = 0.183992005330519
$ws.Range("G3").Value = 155.8700226666667
$ws.Range("H3").Value = 467.610068
$ws.Range("I3").Value = 0.4627663557222626
$ws.Range("J3").Value = 0.4864916976605717
$ws.Range("O3").Value = 0.02891473894686308
$ws.Range("P3").Value = 0.03736411928828315
$ws.Range("Q3").Value = 2.484672074655112
$ws.Range("R3").Value = 22.362048671896
$ws.Range("S3").Value = 0.0133807683691004
$ws.Range("T3").Value = 0.01817733382414898
$ws.Range("G4").Value = 155.8700226666667
$ws.Range("H4").Value = 467.610068
$ws.Range("I4").Value = 0.4627663557222626
$ws.Range("J4").Value = 0.4864916976605717
$ws.Range("M4").Value = 0.3740059999999999
$ws.Range("N4").Value = 0.7480119999999999
$ws.Range("O4").Value = 0.6784086312509182
$ws.Range("P4").Value = 0.5844341432200085
$ws.Range("Q4").Value = 58.29632369746933
$ws.Range("R4").Value = 349.777942184816
$ws.Range("S4").Value = 0.3139446899745157
$ws.Range("T4").Value = 0.2843223585059036
$ws.Range("I5").Value = 0.3897411505765819
$ws.Range("J5").Value = 0.4097225989911443
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1613523333333333
$ws.Range("N5").Value = 0.484057
$ws.Range("O5").Value = 0.2926766298022186
$ws.Range("P5").Value = 0.3782017374917083
$ws.Range("Q5").Value = 21.18128649794733
$ws.Range("R5").Value = 190.631578481526
$ws.Range("S5").Value = 0.114068126445993
$ws.Range("T5").Value = 0.1549577988280692
$ws.Range("I6").Value = 0.3897411505765819
$ws.Range("J6").Value = 0.4097225989911443
$ws.Range("O6").Value = 0.02891473894686308
$ws.Range("P6").Value = 0.03736411928828315
$ws.Range("S6").Value = 0.01126926362577192
$ws.Range("T6").Value = 0.01530892406381052
$ws.Range("I7").Value = 0.3897411505765819
$ws.Range("J7").Value = 0.4097225989911443
$ws.Range("M7").Value = 0.3740059999999999
$ws.Range("N7").Value = 0.7480119999999999
$ws.Range("O7").Value = 0.6784086312509182
$ws.Range("P7").Value = 0.5844341432200085
$ws.Range("Q7").Value = 49.09707888503599
$ws.Range("R7").Value = 294.582473310216
$ws.Range("S7").Value = 0.264403760504817
$ws.Range("T7").Value = 0.2394558760992645
$ws.Range("G8").Value = 0.2461213333333333
$ws.Range("H8").Value = 0.738364
$ws.Range("I8").Value = 0.0007307156985262189
$ws.Range("J8").Value = 0.0007681784042585035
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1613523333333333
$ws.Range("N8").Value = 0.484057
$ws.Range("O8").Value = 0.2926766298022186
$ws.Range("P8").Value = 0.3782017374917083
$ws.Range("Q8").Value = 0.03971225141644445
$ws.Range("R8").Value = 0.357410262748
$ws.Range("S8").Value = 0.0002138634079882278
$ws.Range("T8").Value = 0.0002905264071941739
$ws.Range("G9").Value = 0.2461213333333333
$ws.Range("H9").Value = 0.738364
$ws.Range("I9").Value = 0.0007307156985262189
$ws.Range("J9").Value = 0.0007681784042585035
$ws.Range("O9").Value = 0.02891473894686308
$ws.Range("P9").Value = 0.03736411928828315
$ws.Range("Q9").Value = 0.003923338134222223
$ws.Range("R9").Value = 0.03531004320800001
$ws.Range("S9").Value = [double]"2.112845366726032E-05"
$ws.Range("T9").Value = [double]"2.870230953139772E-05"
$ws.Range("G10").Value = 0.2461213333333333
$ws.Range("H10").Value = 0.738364
$ws.Range("I10").Value = 0.0007307156985262189
$ws.Range("J10").Value = 0.0007681784042585035
$ws.Range("M10").Value = 0.3740059999999999
$ws.Range("N10").Value = 0.7480119999999999
$ws.Range("O10").Value = 0.6784086312509182
$ws.Range("P10").Value = 0.5844341432200085
$ws.Range("Q10").Value = 0.09205085539466665
$ws.Range("R10").Value = 0.5523051323679999
$ws.Range("S10").Value = 0.0004957238368707308
$ws.Range("T10").Value = 0.0004489496875329319
$ws.Range("G11").Value = 49.2786865
$ws.Range("H11").Value = 98.557373
$ws.Range("I11").Value = 0.1463047080910041
$ws.Range("J11").Value = 0.1025370217386683
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1613523333333333
$ws.Range("N11").Value = 0.484057
$ws.Range("O11").Value = 0.2926766298022186
$ws.Range("P11").Value = 0.3782017374917083
$ws.Range("Q11").Value = 7.951231050376834
$ws.Range("R11").Value = 47.707386302261
$ws.Range("S11").Value = 0.04281996888827247
$ws.Range("T11").Value = 0.03877967977878943
$ws.Range("G12").Value = 49.2786865
$ws.Range("H12").Value = 98.557373
$ws.Range("I12").Value = 0.1463047080910041
$ws.Range("J12").Value = 0.1025370217386683
$ws.Range("O12").Value = 0.02891473894686308
$ws.Range("P12").Value = 0.03736411928828315
$ws.Range("Q12").Value = 0.7855351152676667
$ws.Range("R12").Value = 4.713210691606
$ws.Range("S12").Value = 0.00423036244114839
$ws.Range("T12").Value = 0.003831205511708886
$ws.Range("G13").Value = 49.2786865
$ws.Range("H13").Value = 98.557373
$ws.Range("I13").Value = 0.1463047080910041
$ws.Range("J13").Value = 0.1025370217386683
$ws.Range("M13").Value = 0.3740059999999999
$ws.Range("N13").Value = 0.7480119999999999
$ws.Range("O13").Value = 0.6784086312509182
$ws.Range("P13").Value = 0.5844341432200085
$ws.Range("Q13").Value = 18.430524423119
$ws.Range("R13").Value = 73.72209769247598
$ws.Range("S13").Value = 0.09925437676158323
$ws.Range("T13").Value = 0.05992613644817002
$ws.Range("G14").Value = 0.1539513333333333
$ws.Range("H14").Value = 0.461854
$ws.Range("I14").Value = 0.0004570699116250634
$ws.Range("J14").Value = 0.0004805032053572586
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1613523333333333
$ws.Range("N14").Value = 0.484057
$ws.Range("O14").Value = 0.2926766298022186
$ws.Range("P14").Value = 0.3782017374917083
$ws.Range("Q14").Value = 0.02484040685311111
$ws.Range("R14").Value = 0.223563661678
$ws.Range("S14").Value = 0.0001337736813184215
$ws.Range("T14").Value = 0.0001817271471364503
$ws.Range("G15").Value = 0.1539513333333333
$ws.Range("H15").Value = 0.461854
$ws.Range("I15").Value = 0.0004570699116250634
$ws.Range("J15").Value = 0.0004805032053572586
$ws.Range("O15").Value = 0.02891473894686308
$ws.Range("P15").Value = 0.03736411928828315
$ws.Range("Q15").Value = 0.002454086887555556
$ws.Range("R15").Value = 0.022086781988
$ws.Range("S15").Value = [double]"1.321605717510448E-05"
$ws.Range("T15").Value = [double]"1.795357908337102E-05"
$ws.Range("G16").Value = 0.1539513333333333
$ws.Range("H16").Value = 0.461854
$ws.Range("I16").Value = 0.0004570699116250634
$ws.Range("J16").Value = 0.0004805032053572586
$ws.Range("M16").Value = 0.3740059999999999
$ws.Range("N16").Value = 0.7480119999999999
$ws.Range("O16").Value = 0.6784086312509182
$ws.Range("P16").Value = 0.5844341432200085
$ws.Range("Q16").Value = 0.05757872237466666
$ws.Range("R16").Value = 0.3454723342479999
$ws.Range("S16").Value = 0.0003100801731315374
$ws.Range("T16").Value = 0.0002808224791374372
